# Generate Report for Handoff
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": the status text and the two "last generated"
# timestamps that track that transition are refreshed across the
# Overview sheet and each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus  = "Ready for handoff"
$overviewTs = "2016-08-24 21:01:05"
$zhHoTs     = "2016-08-24 21:00:57"
$deHoTs     = "2016-08-24 21:01:05"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $overviewTs

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $zhHoTs

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $deHoTs

# --- Widen the "Status" columns so the longer text fits, matching the
#     auto-fit Excel performs after the cell content changes. -----------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
